$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so values like "0.9989" or
# "0.00000000120" are not auto-converted to numbers by Excel.
$rngFmt = $ws.Range("D2:E51")
$rngFmt.NumberFormat = "@"

$ws.Range("D2").Value = "29.050.09"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.827.00"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "240.51"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "0.6175"
$ws.Range("E6").Value = "  -7.18%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "44.52"
$ws.Range("E8").Value = "  +6.43%  "
$ws.Range("D9").Value = "0.07333"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "0.2913"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "22.66"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "0.07695"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "1.826.98"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "0.6612"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "81.64"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "0.000008935"
$ws.Range("E17").Value = "  +6.41%  "
$ws.Range("D18").Value = "6.023"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "29.038.57"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "2.073.52"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "224.76"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").Value = "12.36"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "7.115"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "1.000"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "159.93"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "8.420"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").Value = "0.1348"
$ws.Range("E28").Value = "  -4.36%  "
$ws.Range("D29").Value = "17.75"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "1.491"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "4.033"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.046"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "1.198"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "0.05270"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").Value = "1.840"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").Value = "1.146"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").Value = "0.7295"
$ws.Range("E37").Value = "  -3.95%  "
$ws.Range("D38").Value = "2.647"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "1.287.25"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.746"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01781"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "6.324"
$ws.Range("E42").Value = "  +6.10%  "
$ws.Range("D43").Value = "0.8993"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").Value = "0.9994"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "102.02"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "1.972.98"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "0.5114"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000120"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "63.65"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.715"
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "0.3960"
$ws.Range("E51").Value = "  -1.90%  "

# Restore default (unstyled) cell style now that the text values are set,
# so the output keeps matching the original styling of these data cells.
$rngFmt.Style = "Normal"
